$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the JDBC-statement flag column: it should hold the literal text "yes",
# not the numeric placeholder 1 that was there before.
$ws.Range("Q1").Value = "yes"

# Leave the cursor where the author ended up after making the edit.
$ws.Range("Q4").Select()
